{"js": "// Update the division problems in the single table of the worksheet.\n// Each cell's text is replaced in place (paragraph range replace) so that\n// the existing run/paragraph formatting (font, size, alignment) is kept\n// untouched and only the visible digits change, matching the source diff.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (rowIndex, colIndex, expectedOldText, newText) \u2013 rows 0,4,8,12,16 are the\n// five populated rows of the practice-problem table; the rows in between\n// are intentionally blank spacer rows used for student work.\nconst replacements = [\n  [0, 0, \"82\u00f72=\", \"99\u00f74=\"],\n  [0, 1, \"24\u00f73=\", \"26\u00f76=\"],\n  [0, 2, \"66\u00f72=\", \"42\u00f73=\"],\n  [0, 3, \"92\u00f75=\", \"73\u00f77=\"],\n  [0, 4, \"99\u00f77=\", \"76\u00f78=\"],\n\n  [4, 0, \"78\u00f73=\", \"61\u00f76=\"],\n  [4, 1, \"65\u00f75=\", \"46\u00f74=\"],\n  [4, 2, \"48\u00f77=\", \"34\u00f76=\"],\n  [4, 3, \"25\u00f76=\", \"47\u00f74=\"],\n  [4, 4, \"77\u00f74=\", \"75\u00f73=\"],\n\n  [8, 0, \"80\u00f72=\", \"48\u00f78=\"],\n  [8, 1, \"45\u00f78=\", \"24\u00f72=\"],\n  [8, 2, \"17\u00f76=\", \"99\u00f78=\"],\n  [8, 3, \"11\u00f76=\", \"43\u00f74=\"],\n  [8, 4, \"30\u00f72=\", \"67\u00f75=\"],\n\n  [12, 0, \"55\u00f78=\", \"69\u00f75=\"],\n  [12, 1, \"67\u00f75=\", \"37\u00f73=\"],\n  [12, 2, \"71\u00f75=\", \"82\u00f78=\"],\n  [12, 3, \"45\u00f73=\", \"60\u00f76=\"],\n  [12, 4, \"69\u00f72=\", \"97\u00f74=\"],\n\n  [16, 0, \"57\u00f76=\", \"51\u00f77=\"],\n  [16, 1, \"78\u00f77=\", \"84\u00f74=\"],\n  [16, 2, \"74\u00f78=\", \"18\u00f79=\"],\n  [16, 3, \"49\u00f77=\", \"39\u00f79=\"],\n  [16, 4, \"45\u00f75=\", \"69\u00f78=\"],\n];\n\n// Load the current text of every target cell first (single sync) so we can\n// confirm we're editing the right cell before writing, without assuming\n// anything about text ordering collisions (some new values equal other old\n// values elsewhere in the table).\nconst cells = replacements.map(([r, c]) => table.getCell(r, c));\ncells.forEach((cell) => cell.body.load(\"text\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, , oldText, newText] = replacements[i];\n  const cell = cells[i];\n  const actual = cell.body.text.replace(/[\\r\\n]+$/, \"\");\n  if (actual !== oldText) {\n    throw new Error(\n      `Unexpected cell text at replacement ${i}: expected \"${oldText}\" but found \"${actual}\"`\n    );\n  }\n  const range = cell.body.paragraphs.getFirst().getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the division problems in the single table of the worksheet.\n# Each cell's Range.Text is replaced in place so the existing run/paragraph\n# formatting (font, size, alignment) stays untouched and only the visible\n# digits change, matching the source diff.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# (rowIndex, colIndex, expectedOldText, newText) using 1-based COM indices.\n# Rows 1,5,9,13,17 are the five populated rows of the practice-problem\n# table; the rows in between are intentionally blank spacer rows used for\n# student work.\n$replacements = @(\n    @(1, 1, \"82\u00f72=\", \"99\u00f74=\"),\n    @(1, 2, \"24\u00f73=\", \"26\u00f76=\"),\n    @(1, 3, \"66\u00f72=\", \"42\u00f73=\"),\n    @(1, 4, \"92\u00f75=\", \"73\u00f77=\"),\n    @(1, 5, \"99\u00f77=\", \"76\u00f78=\"),\n\n    @(5, 1, \"78\u00f73=\", \"61\u00f76=\"),\n    @(5, 2, \"65\u00f75=\", \"46\u00f74=\"),\n    @(5, 3, \"48\u00f77=\", \"34\u00f76=\"),\n    @(5, 4, \"25\u00f76=\", \"47\u00f74=\"),\n    @(5, 5, \"77\u00f74=\", \"75\u00f73=\"),\n\n    @(9, 1, \"80\u00f72=\", \"48\u00f78=\"),\n    @(9, 2, \"45\u00f78=\", \"24\u00f72=\"),\n    @(9, 3, \"17\u00f76=\", \"99\u00f78=\"),\n    @(9, 4, \"11\u00f76=\", \"43\u00f74=\"),\n    @(9, 5, \"30\u00f72=\", \"67\u00f75=\"),\n\n    @(13, 1, \"55\u00f78=\", \"69\u00f75=\"),\n    @(13, 2, \"67\u00f75=\", \"37\u00f73=\"),\n    @(13, 3, \"71\u00f75=\", \"82\u00f78=\"),\n    @(13, 4, \"45\u00f73=\", \"60\u00f76=\"),\n    @(13, 5, \"69\u00f72=\", \"97\u00f74=\"),\n\n    @(17, 1, \"57\u00f76=\", \"51\u00f77=\"),\n    @(17, 2, \"78\u00f77=\", \"84\u00f74=\"),\n    @(17, 3, \"74\u00f78=\", \"18\u00f79=\"),\n    @(17, 4, \"49\u00f77=\", \"39\u00f79=\"),\n    @(17, 5, \"45\u00f75=\", \"69\u00f78=\")\n)\n\nforeach ($entry in $replacements) {\n    $row = $entry[0]\n    $col = $entry[1]\n    $oldText = $entry[2]\n    $newText = $entry[3]\n\n    $cell = $t.Cell($row, $col)\n    $range = $cell.Range\n    $actual = $range.Text -replace \"[\\r\\a]+$\", \"\"\n\n    if ($actual -ne $oldText) {\n        throw \"Unexpected cell text at row $row col ${col}: expected '$oldText' but found '$actual'\"\n    }\n\n    $range.Text = $newText\n}\n"}
